{"js": "// Highlight the first part of three specific run texts (adding SQL Search\n// functions section emphasis) by searching for the literal run-opening\n// phrase and setting font.highlightColor = \"Yellow\" on just that portion.\n// Word (and this engine) splits the underlying run in two so only the\n// matched span carries the new <w:highlight w:val=\"yellow\"/>, leaving the\n// remainder of the original run's text/formatting untouched.\n\nconst body = context.document.body;\n\nconst targets = [\n  \"Nessa tela ser\u00e1 exibido todas as informa\u00e7\u00f5es do produto\",\n  \"Aqui ser\u00e3o exibidos em uma tabela todos os clientes\",\n  \"Aqui ser\u00e1 exibido todas as informa\u00e7\u00f5es do cliente\"\n];\n\nfor (const needle of targets) {\n  const results = body.search(needle, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + needle);\n  }\n\n  const hit = results.items[0];\n  hit.font.highlightColor = \"Yellow\";\n}\n\nawait context.sync();\n", "ps1": "# Highlight the first part of three specific run texts (adding SQL Search\n# functions section emphasis) by finding the literal run-opening phrase\n# and setting Font.HighlightColorIndex = wdYellow (7) on just that span.\n# Word splits the underlying run in two so only the matched span carries\n# the new <w:highlight w:val=\"yellow\"/>, leaving the remainder of the\n# original run's text/formatting untouched.\n\n$d = $word.ActiveDocument\n\n$targets = @(\n  \"Nessa tela ser\u00e1 exibido todas as informa\u00e7\u00f5es do produto\",\n  \"Aqui ser\u00e3o exibidos em uma tabela todos os clientes\",\n  \"Aqui ser\u00e1 exibido todas as informa\u00e7\u00f5es do cliente\"\n)\n\nforeach ($needle in $targets) {\n  $rng = $d.Content\n  $found = $rng.Find.Execute($needle)\n  if (-not $found) {\n    throw \"Text not found: $needle\"\n  }\n  $rng.Font.HighlightColorIndex = 7\n}\n"}
